# Add "scroll to specific sub item" logic:
# - Row 17 gets a new note in column B describing the completion date
#   for the newly added feature.
# - Selection moves to the newly edited cell (B17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B17").Value = "2022/3/29完成"

[void]$ws.Range("B17").Select()
